# Applies the "error solve ifrs list" correction to the IFRS financial
# data table on the active sheet. Rows 2-6 get their per-period financial
# figures rescaled/corrected (column D onward); rows 7-9 had their
# (now-invalid) financial figures removed entirely, leaving only the
# label columns (A/B/C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = 4804
$ws.Range("E2").Value = 25
$ws.Range("F2").Value = 25
$ws.Range("G2").Value = -13
$ws.Range("H2").Value = -17
$ws.Range("I2").Value = -23
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 5240
$ws.Range("L2").Value = 2557
$ws.Range("M2").Value = 2683
$ws.Range("N2").Value = 2686
$ws.Range("O2").Value = -3
$ws.Range("P2").Value = 123
$ws.Range("Q2").Value = 141
$ws.Range("R2").Value = -45
$ws.Range("S2").Value = -44
$ws.Range("T2").Value = 105
$ws.Range("U2").Value = 36
$ws.Range("V2").Value = 1511
$ws.Range("W2").Value = 0.52
$ws.Range("X2").Value = -0.36
$ws.Range("Y2").Value = -0.85
$ws.Range("Z2").Value = -0.33
$ws.Range("AA2").Value = 95.3
$ws.Range("AB2").Value = 2137.18
$ws.Range("AC2").Value = -191
$ws.Range("AD2").Value = -35.93
$ws.Range("AE2").Value = 22070
$ws.Range("AF2").Value = 0.31
$ws.Range("AG2").Value = 146
$ws.Range("AH2").Value = 2.12
$ws.Range("AI2").Value = -78.11
$ws.Range("AJ2").Value = 12607989

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = 5514
$ws.Range("E3").Value = 116
$ws.Range("F3").Value = 116
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 42
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 5636
$ws.Range("L3").Value = 2906
$ws.Range("M3").Value = 2729
$ws.Range("N3").Value = 2682
$ws.Range("O3").Value = 48
$ws.Range("P3").Value = 124
$ws.Range("Q3").Value = 248
$ws.Range("R3").Value = -309
$ws.Range("S3").Value = 173
$ws.Range("T3").Value = 125
$ws.Range("U3").Value = 123
$ws.Range("V3").Value = 1769
$ws.Range("W3").Value = 2.1
$ws.Range("X3").Value = 0.9
$ws.Range("Y3").Value = 1.58
$ws.Range("Z3").Value = 0.92
$ws.Range("AA3").Value = 106.48
$ws.Range("AB3").Value = 2102.5
$ws.Range("AC3").Value = 336
$ws.Range("AD3").Value = 35.83
$ws.Range("AE3").Value = 22038
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 165
$ws.Range("AH3").Value = 1.37
$ws.Range("AI3").Value = 47.21
$ws.Range("AJ3").Value = 12607989

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = 5512
$ws.Range("E4").Value = 139
$ws.Range("F4").Value = 139
$ws.Range("G4").Value = 117
$ws.Range("H4").Value = 80
$ws.Range("I4").Value = 86
$ws.Range("J4").Value = -6
$ws.Range("K4").Value = 5779
$ws.Range("L4").Value = 2954
$ws.Range("M4").Value = 2825
$ws.Range("N4").Value = 2776
$ws.Range("O4").Value = 49
$ws.Range("P4").Value = 126
$ws.Range("Q4").Value = 314
$ws.Range("R4").Value = -218
$ws.Range("S4").Value = -34
$ws.Range("T4").Value = 193
$ws.Range("U4").Value = 121
$ws.Range("V4").Value = 1797
$ws.Range("W4").Value = 2.52
$ws.Range("X4").Value = 1.46
$ws.Range("Y4").Value = 3.16
$ws.Range("Z4").Value = 1.41
$ws.Range("AA4").Value = 104.54
$ws.Range("AB4").Value = 2129.98
$ws.Range("AC4").Value = 683
$ws.Range("AD4").Value = 15.36
$ws.Range("AE4").Value = 22819
$ws.Range("AF4").Value = 0.46
$ws.Range("AG4").Value = 205
$ws.Range("AH4").Value = 1.95
$ws.Range("AI4").Value = 28.94
$ws.Range("AJ4").Value = 12607989

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = 6218
$ws.Range("E5").Value = -8
$ws.Range("F5").Value = -8
$ws.Range("G5").Value = -68
$ws.Range("H5").Value = -94
$ws.Range("I5").Value = -71
$ws.Range("J5").Value = -6
$ws.Range("K5").Value = 6153
$ws.Range("L5").Value = 3502
$ws.Range("M5").Value = 2650
$ws.Range("N5").Value = 2631
$ws.Range("O5").Value = 36
$ws.Range("P5").Value = 126
$ws.Range("Q5").Value = 44
$ws.Range("R5").Value = -278
$ws.Range("S5").Value = 307
$ws.Range("T5").Value = 251
$ws.Range("U5").Value = -206
$ws.Range("V5").Value = 2143
$ws.Range("W5").Value = -0.13
$ws.Range("X5").Value = -1.52
$ws.Range("Y5").Value = -2.62
$ws.Range("Z5").Value = -1.58
$ws.Range("AA5").Value = 132.15
$ws.Range("AB5").Value = 2050.67
$ws.Range("AC5").Value = -563
$ws.Range("AD5").Value = -13.68
$ws.Range("AE5").Value = 21623
$ws.Range("AF5").Value = 0.36
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 12607989

# --- Row 6 (no J6/O6 in either version) -----------------------------------
$ws.Range("D6").Value = 5660
$ws.Range("E6").Value = -71
$ws.Range("F6").Value = -71
$ws.Range("G6").Value = -99
$ws.Range("H6").Value = -92
$ws.Range("I6").Value = -79
$ws.Range("K6").Value = 6253
$ws.Range("L6").Value = 3702
$ws.Range("M6").Value = 2551
$ws.Range("N6").Value = 2544
$ws.Range("P6").Value = 126
$ws.Range("Q6").Value = 239
$ws.Range("R6").Value = -587
$ws.Range("S6").Value = 394
$ws.Range("T6").Value = 498
$ws.Range("U6").Value = -259
$ws.Range("V6").Value = 2515
$ws.Range("W6").Value = -1.26
$ws.Range("X6").Value = -1.62
$ws.Range("Y6").Value = -3.04
$ws.Range("Z6").Value = -1.48
$ws.Range("AA6").Value = 145.13
$ws.Range("AB6").Value = 1983.9
$ws.Range("AC6").Value = -623
$ws.Range("AD6").Value = -8.49
$ws.Range("AE6").Value = 20912
$ws.Range("AF6").Value = 0.25
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 1.32
$ws.Range("AI6").Value = -10.84
$ws.Range("AJ6").Value = 12607989

# --- Rows 7-9: drop the (erroneous) financial figures entirely, keep -----
# --- only the row number / period / label columns (A/B/C). ---------------
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
